$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.066.64'
$ws.Range("E2").Value = '  +0.21%  '
$ws.Range("D3").Value = '2.300.10'
$ws.Range("E3").Value = '  +0.26%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '300.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '97.83'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.17%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.521'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.68%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").Value = '  +1.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.25'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.33%  '
$ws.Range("E11").Value = '  +0.73%  '
$ws.Range("E12").Value = '  +0.59%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '17.77'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.23%  '
$ws.Range("E14").Value = '  -0.53%  '
$ws.Range("D15").Value = '2.658.69'
$ws.Range("E15").Value = '  +0.11%  '
$ws.Range("D16").Value = '2.253.17'
$ws.Range("E16").Value = '  +2.17%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.790'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.12%  '
$ws.Range("D18").Value = '42.948.11'
$ws.Range("E18").Value = '  +0.16%  '
$ws.Range("E19").Value = '  +4.10%  '
$ws.Range("D20").Value = '0.0₃0912'
$ws.Range("E20").Value = '  +1.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.13'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.73%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '237.95'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.99%  '
$ws.Range("E24").Value = '  -1.00%  '
$ws.Range("E25").Value = '  -0.46%  '
$ws.Range("E26").Value = '  -0.45%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '24.95'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.05%  '
$ws.Range("E29").Value = '  -13.12%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.16'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.35%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '163.46'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.09%  '
$ws.Range("E32").Value = '  -3.87%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  +2.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '18.14'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.71%  '
$ws.Range("E36").Value = '  +2.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.42'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0697'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.85%  '
$ws.Range("E39").Value = '  +0.65%  '
$ws.Range("E40").Value = '  -0.10%  '
$ws.Range("E41").Value = '  +1.96%  '
$ws.Range("E42").Value = '  -1.38%  '
$ws.Range("D43").Value = '2.013.65'
$ws.Range("E43").Value = '  +2.49%  '
$ws.Range("E44").Value = '  -1.39%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.19'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.64%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.32'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.43%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '17.39'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.67%  '
$ws.Range("E48").Value = '  -1.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.39'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.74%  '
$ws.Range("D50").Value = '2.530.25'
$ws.Range("E50").Value = '  +0.41%  '
$ws.Range("E51").Value = '  -0.39%  '
